$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.217.96"
$ws.Range("E2").Value = "  +5.76%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.32"
$ws.Range("E3").Value = "  +2.25%  "

$ws.Range("E4").Value = "  -0.49%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.83"
$ws.Range("E5").Value = "  +4.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5185"
$ws.Range("E7").Value = "  +2.07%  "

$ws.Range("E8").Value = "  +3.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08501"
$ws.Range("E9").Value = "  +1.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.127"
$ws.Range("E10").Value = "  +1.93%  "

$ws.Range("E11").Value = "  +1.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.41"
$ws.Range("E12").Value = "  +15.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.449"
$ws.Range("E13").Value = "  +4.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.918.51"
$ws.Range("E14").Value = "  +2.75%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.403"
$ws.Range("E15").Value = "  +1.94%  "

$ws.Range("E16").Value = "  -0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.34"
$ws.Range("E17").Value = "  +2.14%  "

$ws.Range("E18").Value = "  +1.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06696"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.51"
$ws.Range("E20").Value = "  +4.86%  "

$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.019"
$ws.Range("E22").Value = "  +1.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.231.63"
$ws.Range("E23").Value = "  +5.71%  "

$ws.Range("E24").Value = "  +2.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.231"
$ws.Range("E25").Value = "  +1.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.141.56"
$ws.Range("E26").Value = "  +2.72%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.60"
$ws.Range("E27").Value = "  +2.93%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.33"
$ws.Range("E28").Value = "  +3.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.420"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.82"
$ws.Range("E30").Value = "  +2.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.098"
$ws.Range("E31").Value = "  +5.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1067"
$ws.Range("E32").Value = "  +2.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.017"
$ws.Range("E33").Value = "  +4.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.645"
$ws.Range("E34").Value = "  +0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02494"
$ws.Range("E35").Value = "  +1.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06580"
$ws.Range("E36").Value = "  +0.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2214"
$ws.Range("E37").Value = "  +2.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.187"
$ws.Range("E38").Value = "  +3.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.228"
$ws.Range("E39").Value = "  +3.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.92"
$ws.Range("E40").Value = "  +6.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.828"
$ws.Range("E41").Value = "  -2.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6530"
$ws.Range("E42").Value = "  +2.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.241"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6155"
$ws.Range("E44").Value = "  +2.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.36"
$ws.Range("E45").Value = "  +2.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.740"
$ws.Range("E46").Value = "  +1.74%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.080"
$ws.Range("E47").Value = "  +3.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.246"
$ws.Range("E48").Value = "  +2.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.22"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("E50").Value = "  -1.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.56"
$ws.Range("E51").Value = "  +4.31%  "
